$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$shape = $s.Shapes.Item(2)
$shape.TextFrame.TextRange.Paragraphs(4).Text = "User"
